# Apply the updated cryptocurrency price/volume data to the worksheet.
# Numeric-looking text values (e.g. "0.9998") are prefixed with a leading
# apostrophe so Excel stores them as text instead of converting them to
# numbers (matching the original inline-string/text cell content).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.039.60'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '1.886.64'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('D4').Value = '''0.9998'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '''0.7374'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = '''0.9993'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  +2.54%  '
$ws.Range('D9').Value = '''0.07186'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').Value = '''24.77'
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('D11').Value = '''0.08327'
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('D12').Value = '''0.7572'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''5.395'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.851.81'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = '''93.13'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').Value = '''6.149'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '30.032.34'
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').Value = '''250.17'
$ws.Range('E18').Value = '  +4.38%  '
$ws.Range('D19').Value = '''13.56'
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').Value = '''0.000007855'
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.145.24'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '''0.9990'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').Value = '''7.889'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '''0.9995'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').Value = '''0.1559'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').Value = '''9.276'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Value = '''163.24'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('D29').Value = '''2.048'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('D30').Value = '''1.477'
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('D31').Value = '''4.565'
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').Value = '''4.204'
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('D34').Value = '''0.05339'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').Value = '''1.250'
$ws.Range('E35').Value = '  +2.03%  '
$ws.Range('D36').Value = '''0.7686'
$ws.Range('E36').Value = '  +3.61%  '
$ws.Range('D37').Value = '''1.000'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').Value = '''2.758'
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').Value = '''0.4572'
$ws.Range('E41').Value = '  +3.46%  '
$ws.Range('D42').Value = '''6.049'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').Value = '1.091.09'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').Value = '''72.27'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('D45').Value = '''0.8715'
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('D46').Value = '''104.66'
$ws.Range('E46').Value = '  +2.79%  '
$ws.Range('D47').Value = '''0.9999'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').Value = '''1.856'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('D49').Value = '''7.582'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D50').Value = '''9.548'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').Value = '2.019.22'
$ws.Range('E51').Value = '  -0.68%  '
